$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells (row 1): Wins / Losses / Ties, using the same
# formatting (bold, bordered, centered) as the existing header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record columns for every player row (2-52): Wins=84, Losses=78, Ties=0
$lastRow = 52
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 78
    $ws.Cells.Item($r, 32).Value = 0
}
